# Fruta / hortaliza, semanal
# Update the weekly price records: dates and volume/price figures for rows 2-10
# get reshuffled across the records (columns D, M, N, O, P, S), matching the
# new weekly snapshot while leaving all other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$data = @{
    2  = @(44798, 80,  21000, 22000, 21500, 1075)
    3  = @(44761, 100, 20000, 21000, 20500, 1025)
    4  = @(44893, 80,  21000, 22000, 21625, 1081)
    5  = @(44320, 80,  16000, 17000, 16500, 825)
    6  = @(44533, 100, 16000, 17000, 16500, 825)
    7  = @(44357, 100, 14000, 15000, 14500, 725)
    8  = @(44708, 80,  20000, 21000, 20500, 1025)
    9  = @(44890, 80,  20000, 23000, 22250, 1112)
    10 = @(44792, 100, 21000, 22000, 21500, 1075)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value  = $vals[0]  # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]  # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]  # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]  # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]  # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals[5]  # S: Precio $/Kg
}
